$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77 - Myakumyaku (by Memeta)
$ws.Range("A77").Value = "Myakumyaku"
$ws.Range("B77").Value = "Memeta"
$ws.Range("C77").Value = "01haomingHat"
$ws.Range("D77").Value = "myakumyaku.png"
$ws.Range("F77").Value = "myakumyaku_climb.png"

# Row 78 - Mosamosa (by Memeta)
$ws.Range("A78").Value = "Mosamosa"
$ws.Range("B78").Value = "Memeta"
$ws.Range("C78").Value = "01haomingHat"
$ws.Range("D78").Value = "mosamosa_adaptive.png"
$ws.Range("F78").Value = "mosamosa_adaptive_climb.png"

# Row 79 - Swan (by Enoki)
$ws.Range("A79").Value = "Swan"
$ws.Range("B79").Value = "Enoki"
$ws.Range("C79").Value = "01haomingHat"
$ws.Range("D79").Value = "swan_adaptive.png"
$ws.Range("F79").Value = "swan_climb_adaptive.png"

# Match the author's final selection/scroll position in the saved view
$ws.Range("A78").Select()
